$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph by its text.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs($targetIndex)

    # Insert a new empty paragraph right after it.
    $target.Range.InsertParagraphAfter()

    # Fill the new paragraph with the bullet-list entry.
    $newPara = $d.Paragraphs($targetIndex + 1)
    $newPara.Range.Text = "7455355 - Robson da Silva Rocha"
    $newPara.Style = "ListBullet"
}
